# Adicionada a questão 2 ao relatório
# Updates the "x" table (C20:I26) and the "t" table (C28:I34) on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- "x" range (C20:I26): every cell flips from 0 to 1 ---
for ($r = 20; $r -le 26; $r++) {
    for ($c = 3; $c -le 9; $c++) {
        $ws.Cells.Item($r, $c).Value = 1
    }
}

# --- "t" range (C28:I34): updated cumulative totals ---
$tVals = @(
    @(0, 12, 24, 36, 46, 54, 65),
    @(18, 30, 40, 51, 61, 72, 82),
    @(36, 46, 58, 68, 78, 90, 99),
    @(54, 64, 76, 86, 96, 108, 116),
    @(73, 81, 93, 103, 115, 124, 133),
    @(90, 100, 110, 122, 132, 142, 151),
    @(105, 117, 128, 140, 150, 160, 170)
)

for ($i = 0; $i -lt 7; $i++) {
    $row = 28 + $i
    $rowVals = $tVals[$i]
    for ($j = 0; $j -lt 7; $j++) {
        $col = 3 + $j
        $ws.Cells.Item($row, $col).Value = $rowVals[$j]
    }
}
